# Fixed update to excel issue
$wb = $excel.ActiveWorkbook

# 1) Rename header on "Weekly Quantity" sheet (sheet1): B1 "Requested quantity" -> "Weekly_PO_Qty"
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# 2) Rename header on "Monthly Trend" sheet (sheet2): B1 "Requested quantity" -> "Monthly_PO_Qty"
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 3) Add a new "PO Forecast" worksheet at the end of the workbook
$wsForecast = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match header formatting (bold, bordered, centered) used by the other sheets
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = @(
    @(45011.99999999999, 12, -95.07802081835737, 117.0947587444436),
    @(45032.99999999999, 19, -88.98047780663504, 127.6315255694921),
    @(45039.99999999999, 21, -87.82409017190683, 123.5826715889963),
    @(45046.99999999999, 24, -82.03612493993751, 134.3960670876401),
    @(45053.99999999999, 26, -86.98296034735179, 131.3244269105252),
    @(45060.99999999999, 29, -81.55090875442615, 140.7426962846044),
    @(45067.99999999999, 31, -82.3836736894545, 138.9917278749052),
    @(45074.99999999999, 33, -82.38464143826278, 149.1768874825638),
    @(45081.99999999999, 36, -69.47505518077104, 141.0319952274764),
    @(45088.99999999999, 38, -69.828265503148, 150.7593137841582),
    @(45095.99999999999, 40, -60.59016919574522, 153.6536499041881),
    @(45102.99999999999, 43, -66.17203152836754, 151.0165414220361),
    @(45109.99999999999, 45, -62.78750726200983, 153.1684991736801),
    @(45123.99999999999, 50, -55.89125501618859, 159.6156917399529),
    @(45144.99999999999, 57, -56.39752921077044, 166.3927578933361),
    @(45151.99999999999, 60, -55.88866293513821, 168.9217736197962),
    @(45158.99999999999, 62, -46.5053803502475, 180.0712557157354),
    @(45200.99999999999, 76, -30.41380688615384, 181.9912131334709),
    @(45207.99999999999, 79, -31.35548429469441, 180.3761873356731),
    @(45214.99999999999, 81, -23.78364639319732, 190.796358781898),
    @(45221.99999999999, 83, -32.83064438173392, 184.2264939041864),
    @(45242.99999999999, 91, -15.82661340013209, 199.525294245394),
    @(45256.99999999999, 95, -20.35294815406943, 203.3472617502535),
    @(45270.99999999999, 100, -10.40054505655687, 202.80997045991),
    @(45298.99999999999, 110, 3.465453162240133, 217.7337870975105),
    @(45396.99999999999, 143, 29.8190776729299, 246.4094416617936),
    @(45403.99999999999, 146, 38.58709537113348, 254.9361945887452),
    @(45410.99999999999, 148, 35.14512060714496, 246.2505517786917),
    @(45417.99999999999, 150, 40.19921216108252, 254.0167398522296),
    @(45424.99999999999, 153, 43.42957511116241, 257.3520172758126),
    @(45431.99999999999, 155, 46.0684089906855, 265.6036048195201),
    @(45445.99999999999, 160, 65.68859895265648, 267.8163601189343),
    @(45452.99999999999, 162, 55.70590965792977, 270.6500729267889),
    @(45459.99999999999, 165, 64.26609839629606, 269.8339920854724),
    @(45466.99999999999, 167, 66.26628601818443, 270.6660341494683),
    @(45473.99999999999, 170, 60.41073029072845, 279.4119246150148),
    @(45480.99999999999, 172, 66.69766540736575, 287.5221214801343),
    @(45487.99999999999, 174, 67.02973592643522, 284.4969671409735),
    @(45494.99999999999, 177, 69.26669587781858, 289.0214798361019),
    @(45508.99999999999, 181, 66.91592910561332, 281.7055490089031),
    @(45529.99999999999, 189, 73.12387179540681, 295.728833138255),
    @(45564.99999999999, 201, 90.09225281320452, 310.6273157332001),
    @(45578.99999999999, 205, 99.8866993121629, 324.8057479877801),
    @(45599.99999999999, 213, 108.2397875255479, 321.8371322134994),
    @(45606.99999999999, 215, 109.1460020610609, 315.0838746996615),
    @(45613.99999999999, 217, 100.1913802277943, 323.4590636240389),
    @(45620.99999999999, 220, 113.986211917436, 333.7001148411035),
    @(45627.99999999999, 222, 119.866164722228, 324.8991574426324),
    @(45634.99999999999, 224, 116.8241572503269, 336.0006252035013),
    @(45641.99999999999, 227, 116.6054707499153, 335.8270929149903),
    @(45648.99999999999, 229, 120.3519461112967, 341.6708798576693),
    @(45655.99999999999, 232, 119.3204543759949, 335.5449566735389)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Apply the same date/time number format used for the "ds" column on the other sheets
$wsForecast.Range("A2:A53").NumberFormat = "YYYY-MM-DD HH:MM:SS"
